# hours update and TAR update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new status-report rows (17 and 18) for 1/26/2010
# Force the date text to be stored as a literal string (matching the rest of
# the "Date" column, which stores dates typed as text) instead of having it
# auto-parsed into a date serial number: temporarily mark the cell as Text
# so the value is taken verbatim, then restore the date display format.
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "1/26/2010"
$ws.Cells.Item(17, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(17, 2).Value = 2
$ws.Cells.Item(17, 3).Value = "Group Meeting"

$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "1/26/2010"
$ws.Cells.Item(18, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = "Weekly Meeting"

# Update selection to reflect the new last empty row
$ws.Range("A19").Select()
